# Appends two activity-log rows under the existing header row on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "1d59666de757498f"
$ws.Range("B2").Value = "Ali"
$ws.Range("C2").Value = "02-12-2024 18:59:21"
$ws.Range("D2").Value = "02-12-2024 18:59:27"
$ws.Range("E2").Value = "0:00:06"

# Row 3
$ws.Range("A3").Value = "00cee69c916a304e"
$ws.Range("B3").Value = "Ali"
$ws.Range("C3").Value = "02-12-2024 18:59:31"
$ws.Range("D3").Value = "02-12-2024 18:59:32"
$ws.Range("E3").Value = "0:00:01"

# Reset the active selection back to A1 (matches the saved view state).
$ws.Range("A1").Select()
